$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 3 new rows before the current row 35 (FIELDPOLY row) ---
$ws.Range("A35:A37").EntireRow.Insert()

# Give the 3 inserted rows (35 blank separator / 36 header / 37 blank
# separator) their own distinct (but still default-looking) formatting,
# mirroring the row-level custom format Excel leaves behind here.
$ws.Rows("35:37").NumberFormat = "General"

# --- New blank separator row 35 (kept blank, just formatted) ---
# --- New header row 36, mirroring row 5 header with an extra "UNIQUE" column ---
$ws.Range("A36").Value = "HEADER"
$ws.Range("B36").Value = "NAME"
$ws.Range("C36").Value = "DESCRIPTION"
$ws.Range("D36").Value = "ASSOCIATION"
$ws.Range("E36").Value = "UNIT"
$ws.Range("F36").Value = "DATA-TYPE"
$ws.Range("G36").Value = "UNIQUE"

# --- New blank separator row 37 (kept blank, just formatted) ---

# --- Add new "UNIQUE" column (G) values to the 3 properties rows that now sit at 38-40 ---
$ws.Range("G38").Value = "NONE"
$ws.Range("G39").Value = "UNIVERSAL"
$ws.Range("G40").Value = "TYPE"

# --- Selection bookkeeping to match the saved view state ---
$ws.Range("B41").Select()
